$d = $word.ActiveDocument

# Change 1: The first bullet's text was split across three runs
# ("Next Track" / onClick / should update...) with proofErr spell-check
# markers around "onClick". Collapse it back into a single run with no
# proofErr markers by doing a Find & Replace over the whole sentence.
$range = $d.Content
$range.Find.Execute(
    "“Next Track” onClick should update the currently playing song information",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "“Next Track” onClick should update the currently playing song information",
    2
)

# Change 2: The last (previously empty) bullet paragraph gets new text.
$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
$lastPara.Range.Text = "Bank Song button on the side bar should update the displayed song bank if it is the current operation"
$lastPara.Range.Font.Size = 12
$lastPara.Range.Font.SizeBi = 12
